$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 4 data rows (rows 2-5, years 1985-1988), shifting the
# remaining rows up so the data now starts at row 2 with the former row 6
# and ends at row 38 with the former row 42.
$ws.Range("A2:E5").Delete(-4162)
